# Daily TGP refresh: shift existing "today" rows into "yesterday" slot
# and populate new "today" prices (effective date +1 day) per terminal.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(8, 1).Value = 45981
$ws.Cells.Item(8, 4).Value = 172.2
$ws.Cells.Item(8, 5).Value = 162.07
$ws.Cells.Item(8, 6).Value = 172.07
$ws.Cells.Item(8, 7).Value = 162.23

$ws.Cells.Item(9, 1).Value = 45981
$ws.Cells.Item(9, 4).Value = 172.2
$ws.Cells.Item(9, 5).Value = 162.07
$ws.Cells.Item(9, 6).Value = 172.07
$ws.Cells.Item(9, 7).Value = 162.23

$ws.Cells.Item(10, 1).Value = 45981
$ws.Cells.Item(10, 4).Value = 174.24
$ws.Cells.Item(10, 5).Value = 165.06
$ws.Cells.Item(10, 6).Value = 175.06
$ws.Cells.Item(10, 7).Value = 165.55

$ws.Cells.Item(11, 1).Value = 45980
$ws.Cells.Item(11, 4).Value = 171.9
$ws.Cells.Item(11, 5).Value = 161.93
$ws.Cells.Item(11, 6).Value = 171.93
$ws.Cells.Item(11, 7).Value = 162.09

$ws.Cells.Item(12, 1).Value = 45980
$ws.Cells.Item(12, 4).Value = 171.9
$ws.Cells.Item(12, 5).Value = 161.93
$ws.Cells.Item(12, 6).Value = 171.93
$ws.Cells.Item(12, 7).Value = 162.09

$ws.Cells.Item(13, 1).Value = 45980
$ws.Cells.Item(13, 4).Value = 174
$ws.Cells.Item(13, 5).Value = 164.99
$ws.Cells.Item(13, 6).Value = 174.99
$ws.Cells.Item(13, 7).Value = 165.49

$ws.Cells.Item(17, 1).Value = 45981
$ws.Cells.Item(17, 4).Value = 177.55
$ws.Cells.Item(17, 5).Value = 167.83
$ws.Cells.Item(17, 6).Value = 177.83

$ws.Cells.Item(18, 1).Value = 45980
$ws.Cells.Item(18, 4).Value = 177.33
$ws.Cells.Item(18, 5).Value = 167.72
$ws.Cells.Item(18, 6).Value = 177.72

$ws.Cells.Item(22, 1).Value = 45981
$ws.Cells.Item(22, 4).Value = 173.11
$ws.Cells.Item(22, 5).Value = 164.01
$ws.Cells.Item(22, 6).Value = 173.61
$ws.Cells.Item(22, 7).Value = 165.3

$ws.Cells.Item(23, 1).Value = 45981
$ws.Cells.Item(23, 4).Value = 179.02
$ws.Cells.Item(23, 5).Value = 168.55
$ws.Cells.Item(23, 6).Value = 178.55

$ws.Cells.Item(24, 1).Value = 45981
$ws.Cells.Item(24, 4).Value = 178.82
$ws.Cells.Item(24, 5).Value = 168.8
$ws.Cells.Item(24, 6).Value = 178.8

$ws.Cells.Item(25, 1).Value = 45981
$ws.Cells.Item(25, 4).Value = 179.64
$ws.Cells.Item(25, 5).Value = 168.21
$ws.Cells.Item(25, 6).Value = 178.21
$ws.Cells.Item(25, 7).Value = 168.25

$ws.Cells.Item(26, 1).Value = 45981
$ws.Cells.Item(26, 4).Value = 178.34
$ws.Cells.Item(26, 5).Value = 169.78
$ws.Cells.Item(26, 6).Value = 179.78

$ws.Cells.Item(27, 1).Value = 45980
$ws.Cells.Item(27, 4).Value = 172.92
$ws.Cells.Item(27, 5).Value = 163.96
$ws.Cells.Item(27, 6).Value = 173.56
$ws.Cells.Item(27, 7).Value = 165.24

$ws.Cells.Item(28, 1).Value = 45980
$ws.Cells.Item(28, 4).Value = 178.78
$ws.Cells.Item(28, 5).Value = 168.49
$ws.Cells.Item(28, 6).Value = 178.49

$ws.Cells.Item(29, 1).Value = 45980
$ws.Cells.Item(29, 4).Value = 178.58
$ws.Cells.Item(29, 5).Value = 168.73
$ws.Cells.Item(29, 6).Value = 178.73

$ws.Cells.Item(30, 1).Value = 45980
$ws.Cells.Item(30, 4).Value = 179.41
$ws.Cells.Item(30, 5).Value = 168.14
$ws.Cells.Item(30, 6).Value = 178.14
$ws.Cells.Item(30, 7).Value = 168.18

$ws.Cells.Item(31, 1).Value = 45980
$ws.Cells.Item(31, 4).Value = 178.1
$ws.Cells.Item(31, 5).Value = 169.72
$ws.Cells.Item(31, 6).Value = 179.72

$ws.Cells.Item(35, 1).Value = 45981
$ws.Cells.Item(35, 4).Value = 172.49
$ws.Cells.Item(35, 5).Value = 162.05
$ws.Cells.Item(35, 6).Value = 171.05

$ws.Cells.Item(36, 1).Value = 45980
$ws.Cells.Item(36, 4).Value = 172.36
$ws.Cells.Item(36, 5).Value = 161.99
$ws.Cells.Item(36, 6).Value = 170.99

$ws.Cells.Item(40, 1).Value = 45981
$ws.Cells.Item(40, 4).Value = 178.24
$ws.Cells.Item(40, 5).Value = 167.76
$ws.Cells.Item(40, 6).Value = 177.76

$ws.Cells.Item(41, 1).Value = 45981
$ws.Cells.Item(41, 4).Value = 177.94
$ws.Cells.Item(41, 5).Value = 168.18
$ws.Cells.Item(41, 6).Value = 178.18

$ws.Cells.Item(42, 1).Value = 45980
$ws.Cells.Item(42, 4).Value = 177.99
$ws.Cells.Item(42, 5).Value = 167.64
$ws.Cells.Item(42, 6).Value = 177.64

$ws.Cells.Item(43, 1).Value = 45980
$ws.Cells.Item(43, 4).Value = 177.69
$ws.Cells.Item(43, 5).Value = 168.06
$ws.Cells.Item(43, 6).Value = 178.06

$ws.Cells.Item(47, 1).Value = 45981
$ws.Cells.Item(47, 4).Value = 171.88
$ws.Cells.Item(47, 5).Value = 163.43
$ws.Cells.Item(47, 6).Value = 173.43

$ws.Cells.Item(48, 1).Value = 45981
$ws.Cells.Item(48, 4).Value = 171.84
$ws.Cells.Item(48, 5).Value = 163.59
$ws.Cells.Item(48, 6).Value = 173.59

$ws.Cells.Item(49, 1).Value = 45980
$ws.Cells.Item(49, 4).Value = 172.7
$ws.Cells.Item(49, 5).Value = 163.42
$ws.Cells.Item(49, 6).Value = 173.42

$ws.Cells.Item(50, 1).Value = 45980
$ws.Cells.Item(50, 4).Value = 172.66
$ws.Cells.Item(50, 5).Value = 163.58
$ws.Cells.Item(50, 6).Value = 173.58

$ws.Cells.Item(54, 1).Value = 45981
$ws.Cells.Item(54, 4).Value = 188.29
$ws.Cells.Item(54, 5).Value = 178.3
$ws.Cells.Item(54, 6).Value = 188.3

$ws.Cells.Item(55, 1).Value = 45981
$ws.Cells.Item(55, 4).Value = 175.99
$ws.Cells.Item(55, 5).Value = 175.08
$ws.Cells.Item(55, 6).Value = 185.08

$ws.Cells.Item(56, 1).Value = 45981
$ws.Cells.Item(56, 4).Value = 178.39

$ws.Cells.Item(57, 1).Value = 45981
$ws.Cells.Item(57, 4).Value = 177.87
$ws.Cells.Item(57, 5).Value = 169.35

$ws.Cells.Item(58, 1).Value = 45981
$ws.Cells.Item(58, 4).Value = 173.78
$ws.Cells.Item(58, 5).Value = 165.4
$ws.Cells.Item(58, 6).Value = 175.4

$ws.Cells.Item(59, 1).Value = 45981
$ws.Cells.Item(59, 4).Value = 180.52
$ws.Cells.Item(59, 5).Value = 176.27

$ws.Cells.Item(60, 1).Value = 45980
$ws.Cells.Item(60, 4).Value = 188.05
$ws.Cells.Item(60, 5).Value = 178.19
$ws.Cells.Item(60, 6).Value = 188.19

$ws.Cells.Item(61, 1).Value = 45980
$ws.Cells.Item(61, 4).Value = 175.75
$ws.Cells.Item(61, 5).Value = 175.14
$ws.Cells.Item(61, 6).Value = 185.14

$ws.Cells.Item(62, 1).Value = 45980
$ws.Cells.Item(62, 4).Value = 178.26

$ws.Cells.Item(63, 1).Value = 45980
$ws.Cells.Item(63, 4).Value = 177.76
$ws.Cells.Item(63, 5).Value = 169.41

$ws.Cells.Item(64, 1).Value = 45980
$ws.Cells.Item(64, 4).Value = 173.67
$ws.Cells.Item(64, 5).Value = 165.46
$ws.Cells.Item(64, 6).Value = 175.46

$ws.Cells.Item(65, 1).Value = 45980
$ws.Cells.Item(65, 4).Value = 180.29
$ws.Cells.Item(65, 5).Value = 176.19

